# Fill in the missing xG_home / xG_away / goals_home / goals_away values
# for the last six recorded Parma Calcio 1913 matches of 2020 (rows 10-15):
#   10: Genoa      vs Parma Calcio 1913
#   11: Parma      vs Benevento
#   12: AC Milan   vs Parma Calcio 1913
#   13: Parma      vs Cagliari
#   14: Parma      vs Juventus
#   15: Crotone    vs Parma Calcio 1913
#
# These columns hold their values as text (shared strings) in the workbook,
# just like every other row already does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    10 = @{ D = "1.58372";  E = "1.40499";  F = "1"; G = "2" }
    11 = @{ D = "0.501198"; E = "0.340508"; F = "0"; G = "0" }
    12 = @{ D = "1.8515";   E = "0.603434"; F = "2"; G = "2" }
    13 = @{ D = "0.870374"; E = "0.265316"; F = "0"; G = "0" }
    14 = @{ D = "0.926619"; E = "2.80045";  F = "0"; G = "4" }
    15 = @{ D = "0.7999";   E = "1.50226";  F = "2"; G = "1" }
}

# Write column-by-column (D for every row, then E, then F, then G) to match
# how the sheet was originally populated (a column-oriented data export).
# Each value is written as a literal-text formula and then pasted back as a
# plain value, which stores it as a shared-string cell (matching the rest of
# the sheet) without ever creating a number format override on the cell.
foreach ($col in @("D", "E", "F", "G")) {
    foreach ($row in 10..15) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $text = $data[$row][$col]
        $cell.Formula = '="' + $text + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}
$excel.CutCopyMode = $false

Write-Output "done"
